$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Journal")

# --- Fill in the new "Conception" block (rows 32-34) ---
$ws.Range("A32").Value = "Conception"
$ws.Range("C32").Value = 20
$ws.Range("D32").Value = "Mise en place des obstacles"

$ws.Range("A33").Value = "Conception"
$ws.Range("C33").Value = 60
$ws.Range("D33").Value = "Mise en place des hitBox"

$ws.Range("A34").Value = "Conception"
$ws.Range("C34").Value = 30
$ws.Range("D34").Value = "Convention de nommage du code"

# --- Update the active cell / scroll position to match the new view ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 21
$excel.ActiveWindow.ScrollColumn = 4
$ws.Range("D34").Select()
